# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme name "Office")
#   ppt/theme/theme2.xml -> "Integral"     (clrScheme name "Red Violet")
# theme2.xml is the theme actually bound to the (single) slide master that
# every slide/layout in this deck inherits from, so it is the one reachable
# through the PowerPoint object model's colour-scheme APIs.
#
# The target edit swaps the two parts' contents so that the deck's slides
# render with the "Office" colour palette instead of "Integral". We recreate
# that by pushing the Office palette's 12 theme colours onto the reachable
# ThemeColorScheme (keyed off any slide, since they all share the one
# master/theme). Per DrawingML, RGB() values are packed 0xBBGGRR, i.e. the
# byte order is reversed from the usual 0xRRGGBB hex.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
